# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> linked from the (only) Slide Master, currently
#                             the "Integral" template's "Red Violet" colors
#   ppt/theme/theme2.xml  -> linked from the Notes Master, the default
#                             "Office Theme" / "Office" color scheme
#
# The commit swaps the two themes' contents: the Slide Master's theme
# (theme1.xml) ends up carrying the plain default "Office" color scheme
# while the Notes Master's theme (theme2.xml) keeps the "Red Violet" one
# (font scheme + format scheme are identical between the two themes, so
# only the 12-slot color scheme actually changes).
#
# The PowerPoint object model's ThemeColorScheme (exposed off a Slide /
# SlideRange) maps 1:1 onto the <a:clrScheme> children of the Slide
# Master's theme part, in document order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3
#   8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
# Writing .RGB on each slot persists straight into theme1.xml without
# touching anything else (slides, layouts, masters, theme2.xml, ...).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target = default "Office" color scheme (RRGGBB -> the decimal value
# PowerPoint's ColorFormat.RGB getter/setter uses, i.e. R + G*256 + B*65536).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
